$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mistral")

# Apply the A-column (bordered/bold) style to the new rows (8-24) before setting values,
# so that dimension/style info lines up with the rest of the table.
$ws.Range("A2").Copy()
$ws.Range("A8:A24").PasteSpecial(-4122)

# Set cell values for rows 2 through 24 (column A = label, column B = numeric value)
$ws.Range("A2").Value = "preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse"
$ws.Range("B2").Value = 21156.49535999999
$ws.Range("A3").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_lenNone_gblFalse"
$ws.Range("B3").Value = 42281.28768
$ws.Range("A4").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_lenNone_gblFalse"
$ws.Range("B4").Value = 28220.30131200001
$ws.Range("A5").Value = "preds_ns10_ws32_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse"
$ws.Range("B5").Value = 26220.95360000002
$ws.Range("A6").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_lenNone_gblFalse"
$ws.Range("B6").Value = 84484.292608
$ws.Range("A7").Value = "preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse"
$ws.Range("B7").Value = 1261.19936
$ws.Range("A8").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_sum_fused_rerun_lenNone_gblFalse"
$ws.Range("B8").Value = 53752.95488000003
$ws.Range("A9").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse"
$ws.Range("B9").Value = 53752.95488000003
$ws.Range("A10").Value = "preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_qcache_lenNone_gblFalse"
$ws.Range("B10").Value = 1763.704832
$ws.Range("A11").Value = "preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse"
$ws.Range("B11").Value = 0
$ws.Range("A12").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_lenNone_gblFalse"
$ws.Range("B12").Value = 84484.292608
$ws.Range("A13").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse"
$ws.Range("B13").Value = 84562.57535999996
$ws.Range("A14").Value = "preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse"
$ws.Range("B14").Value = 396.816384
$ws.Range("A15").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_lenNone_gblFalse"
$ws.Range("B15").Value = 42281.28768
$ws.Range("A16").Value = "preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_False_type_max_fused_lenNone_gblFalse"
$ws.Range("B16").Value = 21121.073152
$ws.Range("A17").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse"
$ws.Range("B17").Value = 1343.823872
$ws.Range("A18").Value = "preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse"
$ws.Range("B18").Value = 10488.38144
$ws.Range("A19").Value = "preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_lenNone_gblFalse"
$ws.Range("B19").Value = 5255.987199999999
$ws.Range("A20").Value = "preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse"
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = "preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_snapkv_rerun_lenNone_gblFalse"
$ws.Range("B21").Value = 169679.519744
$ws.Range("A22").Value = "preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse"
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = "preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_qcache_lenNone_gblFalse"
$ws.Range("B23").Value = 492.17536
$ws.Range("A24").Value = "preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse"
$ws.Range("B24").Value = 1344.274432

